# --------------------------------------------------------------------------
# "added design changes to documentation"
#
# Appends two new paragraphs right after the bulleted "Resolved (bool) ..."
# list item (end of the "Comments" data-schema bullet list):
#
#   1. A new bulleted item: "deleted (bool) – owner / admin / document
#      owner can toggle." - same ListParagraph / numbered-list (numId 6)
#      formatting as the rest of the bullets in that list.
#   2. An empty ListParagraph that closes off the list (no more numbering,
#      just a left indent) acting as a spacer before the next heading
#      ("General Workflow UML").
# --------------------------------------------------------------------------

$d = $word.ActiveDocument

$origText = "Resolved (bool) – only admins / document owner ID can resolve."
$placeholder = "__NEW_DELETED_BULLET_PLACEHOLDER__"

# 1. Find the paragraph that currently ends the bulleted list and, in one
#    Find/Replace pass, grow it into three paragraphs: the original text,
#    a placeholder paragraph (to be replaced below with the fully-formatted
#    "deleted (bool) ..." bullet) and a trailing empty paragraph.
$rng = $d.Content
$replacement = $origText + "^p" + $placeholder + "^p"
$found = $rng.Find.Execute($origText, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
if (-not $found) {
    throw "Could not find the 'Resolved (bool)' paragraph to anchor the new content"
}

# 2. Re-locate the placeholder paragraph and its following (still empty)
#    sibling paragraph.
$rng2 = $d.Content
$placeholderFound = $rng2.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $placeholderFound) {
    throw "Could not re-locate the placeholder paragraph"
}
$newPara1 = $rng2.Paragraphs(1)
$newPara2 = $newPara1.Next()

# 3. Replace the placeholder paragraph with the finished "deleted (bool) ..."
#    bullet, split across three runs (matching the authored formatting).
$xmlNewPara1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-IL"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-IL"/></w:rPr><w:t>deleted</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-IL"/></w:rPr><w:t xml:space="preserve"> (bool) – </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-IL"/></w:rPr><w:t>owner / admin / document owner can toggle.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara1.Range.InsertXML($xmlNewPara1) | Out-Null

# 4. Replace the trailing empty paragraph with one that keeps the
#    ListParagraph style but drops the numbering and uses a left indent,
#    matching the paragraph that closes out the list.
$xmlNewPara2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="785"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-IL"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara2.Range.InsertXML($xmlNewPara2)

Write-Host "Inserted 'deleted (bool)' bullet and trailing spacer paragraph."
